# Generate Report for Handback
# Updates the "Latest Target File", "Latest Handback DateTime" and
# "Error Detail" columns (J, K, P) for the 6ec5adcc-... row (row 8) on
# both the zh-cn and de-de sheets, turns the (previously empty) I8 cell
# into a hyperlink pointing at the latest handback markdown file (like
# A8 already does), and widens column P (Error Detail) to fit the new,
# longer text.

$wb = $excel.ActiveWorkbook

# Column width helper: Excel's ColumnWidth (character units) is offset
# from the serialized OOXML column width by a constant ~0.8333 so that
# requesting a stored width of 40 requires setting ColumnWidth to this
# value.
$targetColWidth = 39.166666666666664

$handbackFileName = "6ec5adcc-5f68-401a-8bda-a156c3d54cfa.md"
$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8af52a2a979e0993567356433c9c9ba80b56fdc/e2e/6ec5adcc-5f68-401a-8bda-a156c3d54cfa.md"
$currentHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc5c181add62d5d7f832a3423848253540036288/e2e/6ec5adcc-5f68-401a-8bda-a156c3d54cfa.md"
$errorDetail = "The version of handback file is not the latest, current: " + $currentHandbackUrl + ", latest: " + $latestHandbackUrl + "."

# BGR integer matching the workbook's custom HyperLink font color FF6495ED.
$hyperlinkFontColor = 15570276

function Update-LocaleSheet($ws, $targetFileValue, $handbackDateTime) {
    # Widen the Error Detail column (P) to fit the long message.
    $ws.Columns.Item(16).ColumnWidth = $targetColWidth

    # I8: show + link to the latest handback file, same as A8.
    $ws.Range("I8").Value = $handbackFileName
    $null = $ws.Hyperlinks.Add($ws.Range("I8"), $latestHandbackUrl, "", "", $handbackFileName)
    $ws.Range("I8").Font.Underline = 2
    $ws.Range("I8").Font.Color = $hyperlinkFontColor

    # J8: latest target xliff file name.
    $ws.Range("J8").Value = $targetFileValue

    # K8: latest handback datetime (stored as text, matching the other rows).
    $ws.Range("K8").Value = $handbackDateTime

    # P8: error detail message.
    $ws.Range("P8").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn "6ec5adcc-5f68-401a-8bda-a156c3d54cfa.72b98c74150e12e4715945e67fba58a000c61387.zh-cn.xlf" "2016-08-24 06:44:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe "6ec5adcc-5f68-401a-8bda-a156c3d54cfa.72b98c74150e12e4715945e67fba58a000c61387.de-de.xlf" "2016-08-24 06:44:32"
